# Management_information_KIE_2019.xlsx edit
# - Rename fertilizer agent "KAS" -> "CAN" in all "Fertilization/Agent" (column E) cells
#   of the 1st/2nd/3rd nitrogen application rows, across all four treatment blocks.
# - Clear the stray "others"/date/amount/"Bittersalz" data that had been entered on the
#   "2nd application of insecticide" / "1st application of fungicide" rows (columns B-E)
#   across all four treatment blocks — these values did not belong there.
# - Update the saved sheet view (active cell / selection).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Rename "KAS" -> "CAN" (Fertilization Agent column) ---------------------
$kasCells = @("E47", "E48", "E58", "E59", "E69", "E70", "E71", "E80", "E81", "E82")
foreach ($addr in $kasCells) {
    $ws.Range($addr).Value = "CAN"
}

# --- Clear the mistakenly-filled cells ---------------------------------------
$rowPairs = @(50, 61, 72, 83)
foreach ($base in $rowPairs) {
    $ws.Range("B" + $base).ClearContents()
    $ws.Range("C" + $base).ClearContents()
    $ws.Range("D" + $base).ClearContents()
    $ws.Range("E" + $base).ClearContents()

    $next = $base + 1
    $ws.Range("C" + $next).ClearContents()
    $ws.Range("D" + $next).ClearContents()
    $ws.Range("E" + $next).ClearContents()
}

# --- Update view / selection --------------------------------------------------
$ws.Activate()
$win = $excel.ActiveWindow
$win.ScrollRow = 42
$win.ScrollColumn = 1
$ws.Range("E82").Select()
